$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the spell-check ("splitted") proofing split: the original text is
#    spread over 3 runs with <w:proofErr> markers wrapped around "splitted".
#    Deleting that stretch of text and re-inserting it as plain text collapses
#    it back into a single run and drops the proofErr markers, leaving the
#    visible text unchanged.
# ---------------------------------------------------------------------------
$rngSplit = $d.Content
$found = $rngSplit.Find.Execute(
    "Also User is has its own repositories splitted. So if preferences",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitText = $rngSplit.Text
    $rngSplit.Delete()
    $rngSplit.InsertAfter($splitText)
}

# ---------------------------------------------------------------------------
# 2) Mark the picture's run as <w:noProof/> (Word does this automatically for
#    runs that hold inline graphics so the proofer skips them).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $para.Range.NoProofing = 1
    }
}

# ---------------------------------------------------------------------------
# 3) Append the new "Pipeline process" section after the picture paragraph:
#    5 blank paragraphs, a "Pipeline process" heading-like line, a blank
#    line, then three paragraphs of body text.
# ---------------------------------------------------------------------------
function Get-EndRange {
    $r = $d.Content
    $r.Collapse(0)
    return $r
}

# five blank paragraphs
for ($i = 0; $i -lt 5; $i++) {
    $r = Get-EndRange
    $r.InsertParagraphAfter()
}

# "Pipeline process" paragraph
$r = Get-EndRange
$r.InsertParagraphAfter()
$r = Get-EndRange
$r.InsertAfter("Pipeline process")

# blank paragraph
$r = Get-EndRange
$r.InsertParagraphAfter()

# three body paragraphs
$r = Get-EndRange
$r.InsertParagraphAfter()
$r = Get-EndRange
$r.InsertAfter("When the localmachine pushes to git, the runner triggers and runs commands to build, test and start sonar stages.")

$r = Get-EndRange
$r.InsertParagraphAfter()
$r = Get-EndRange
$r.InsertAfter("When everything is executed there will be unit test results and a sonar test report. The test report provides clarity about code smells, duplications test coverage and a lot more, this is all very useful for analyzing code.")

$r = Get-EndRange
$r.InsertParagraphAfter()
$r = Get-EndRange
$r.InsertAfter("It is a very big step on making a project more efficient because tests are automated and code vulnerabilities and bugs will be detected quicker.")

Write-Output "done"
